$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price and volume(1h) values per the scraped data refresh.
# Values are kept as literal text (matching original inlineStr formatting),
# so we force text number-format before assignment and restore the default
# "Normal" style afterwards to avoid leaving stray formatting behind.
$updates = @{
    'D2' = '329.28'
    'E2' = '0.10%'
    'D3' = '44.37'
    'E3' = '0.04%'
    'D4' = '5.521'
    'E4' = '-1.08%'
    'D5' = '0.08110'
    'E5' = '0.24%'
    'D6' = '2.064'
    'E6' = '2.35%'
    'D7' = '0.9738'
    'E7' = '2.11%'
    'D8' = '0.1124'
    'E8' = '-4.39%'
    'D9' = '0.1883'
    'E9' = '1.31%'
    'D10' = '10.15'
    'E10' = '-0.83%'
    'D11' = '0.09948'
    'E11' = '0.90%'
    'D12' = '0.04776'
    'E12' = '2.60%'
    'D13' = '0.1054'
    'E13' = '-1.44%'
    'D14' = '0.001260'
    'E14' = '-2.11%'
    'D15' = '0.04093'
    'E15' = '-3.03%'
    'D16' = '0.005960'
    'E16' = '0.21%'
    'E17' = '-0.92%'
    'D18' = '4.428'
    'E18' = '2.62%'
    'E19' = '3.08%'
    'D20' = '0.3305'
    'D21' = '0.1390'
    'E21' = '-1.35%'
    'D22' = '0.2568'
    'E22' = '2.55%'
    'D23' = '0.001303'
    'E23' = '4.33%'
    'D24' = '0.004388'
    'E24' = '1.62%'
    'E25' = '7.22%'
    'D26' = '0.0003733'
    'E26' = '-6.17%'
    'D38' = '0.02674'
    'E38' = '0.44%'
    'D39' = '0.05642'
    'E39' = '1.34%'
    'D40' = '0.007611'
    'E40' = '1.68%'
    'E41' = '0.24%'
    'D42' = '0.007446'
    'E42' = '-7.84%'
    'D43' = '0.001955'
    'E43' = '-3.12%'
    'D44' = '0.008297'
    'E44' = '-1.29%'
    'D45' = '0.00007076'
    'E45' = '-2.11%'
    'D46' = '0.00000000749'
    'E46' = '-0.25%'
    'D47' = '0.0005794'
    'E47' = '-0.30%'
    'D48' = '0.002516'
    'E48' = '10.73%'
    'D49' = '0.003630'
    'E49' = '-13.49%'
    'D50' = '0.00002096'
    'E50' = '-0.25%'
    'D51' = '0.0001996'
    'E51' = '-0.25%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
